$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.574538
$ws.Cells.Item(2, 8).Value = 1.723614
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.1341725
$ws.Cells.Item(2, 14).Value = 0.268345
$ws.Cells.Item(2, 15).Value = 0.01633512969336317
$ws.Cells.Item(2, 16).Value = 0.01188426112752495
$ws.Cells.Item(2, 17).Value = 0.07708719980499999
$ws.Cells.Item(2, 18).Value = 0.46252319883
$ws.Cells.Item(2, 19).Value = 0.01633512969336317
$ws.Cells.Item(2, 20).Value = 0.01188426112752495

# Row 3: FAPs
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.574538
$ws.Cells.Item(3, 8).Value = 1.723614
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 5.923689
$ws.Cells.Item(3, 14).Value = 17.771067
$ws.Cells.Item(3, 15).Value = 0.7211927040052828
$ws.Cells.Item(3, 16).Value = 0.787031622511101
$ws.Cells.Item(3, 17).Value = 3.403384430682
$ws.Cells.Item(3, 18).Value = 30.630459876138
$ws.Cells.Item(3, 19).Value = 0.7211927040052828
$ws.Cells.Item(3, 20).Value = 0.787031622511101

# Row 4: M1
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.574538
$ws.Cells.Item(4, 8).Value = 1.723614
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.007277666666666668
$ws.Cells.Item(4, 14).Value = 0.021833
$ws.Cells.Item(4, 15).Value = 0.00088603572911786
$ws.Cells.Item(4, 16).Value = 0.0009669234500260939
$ws.Cells.Item(4, 17).Value = 0.004181296051333334
$ws.Cells.Item(4, 18).Value = 0.037631664462
$ws.Cells.Item(4, 19).Value = 0.00088603572911786
$ws.Cells.Item(4, 20).Value = 0.0009669234500260939

# Row 5: M2
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.574538
$ws.Cells.Item(5, 8).Value = 1.723614
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.039371
$ws.Cells.Item(5, 14).Value = 0.118113
$ws.Cells.Item(5, 15).Value = 0.004793310038624915
$ws.Cells.Item(5, 16).Value = 0.005230899530661476
$ws.Cells.Item(5, 17).Value = 0.022620135598
$ws.Cells.Item(5, 18).Value = 0.203581220382
$ws.Cells.Item(5, 19).Value = 0.004793310038624915
$ws.Cells.Item(5, 20).Value = 0.005230899530661476

# Row 6: Neutro
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Wnt2"
$ws.Cells.Item(6, 3).Value = "Fzd2"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.574538
$ws.Cells.Item(6, 8).Value = 1.723614
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.1820473333333333
$ws.Cells.Item(6, 14).Value = 0.5461419999999999
$ws.Cells.Item(6, 15).Value = 0.02216375785150397
$ws.Cells.Item(6, 16).Value = 0.02418712530775207
$ws.Cells.Item(6, 17).Value = 0.1045931107986666
$ws.Cells.Item(6, 18).Value = 0.9413379971879998
$ws.Cells.Item(6, 19).Value = 0.02216375785150397
$ws.Cells.Item(6, 20).Value = 0.02418712530775207

# Row 7: sCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Wnt2"
$ws.Cells.Item(7, 3).Value = "Fzd2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.574538
$ws.Cells.Item(7, 8).Value = 1.723614
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(7, 10).Value = 1
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.927182
$ws.Cells.Item(7, 14).Value = 3.854364
$ws.Cells.Item(7, 15).Value = 0.2346290626821072
$ws.Cells.Item(7, 16).Value = 0.1706991680729343
$ws.Cells.Item(7, 17).Value = 1.107239291916
$ws.Cells.Item(7, 18).Value = 6.643435751496
$ws.Cells.Item(7, 19).Value = 0.2346290626821072
$ws.Cells.Item(7, 20).Value = 0.1706991680729343
